$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12 corresponds to the "WeaponDamage" event (weapon damage tick noise).
# Mark it as Completed and add a note about the sound needing to be crisper.
$ws.Range("E12").Value = "Completed"
$ws.Range("F12").Value = "May need to be a bit more crisp sounding"

# Update the active selection to match the authored edit.
$ws.Range("F12").Select()
